$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2026-01-07 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-08 Thursday", 2) | Out-Null

# Update the 20x5 answer table, cell by cell, in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "0+75=75",
    "75-53=22",
    "73-48=25",
    "76-14=62",
    "1+4=5",
    "47-22=25",
    "20+18=38",
    "22+44=66",
    "64-52=12",
    "96-4=92",
    "21+0=21",
    "64-48=16",
    "25-17=8",
    "43+12=55",
    "8-7=1",
    "74-43=31",
    "68-24=44",
    "23-0=23",
    "10+19=29",
    "33-8=25",
    "24+4=28",
    "26+32=58",
    "22+20=42",
    "91-73=18",
    "52-15=37",
    "91-66=25",
    "40-35=5",
    "91-36=55",
    "64-23=41",
    "59+19=78",
    "27+45=72",
    "49+11=60",
    "90-51=39",
    "9+68=77",
    "65-27=38",
    "71-4=67",
    "1+73=74",
    "14+13=27",
    "62-41=21",
    "92-25=67",
    "92+2=94",
    "82-71=11",
    "26+72=98",
    "41-27=14",
    "86+5=91",
    "19+57=76",
    "40+31=71",
    "41-16=25",
    "20+53=73",
    "70-23=47",
    "96-0=96",
    "42+40=82",
    "86-23=63",
    "95-64=31",
    "14+11=25",
    "18+73=91",
    "59+37=96",
    "93-89=4",
    "86-83=3",
    "37+6=43",
    "74+23=97",
    "16+52=68",
    "21+52=73",
    "3+30=33",
    "92-53=39",
    "5+91=96",
    "78-11=67",
    "1+22=23",
    "63+20=83",
    "25+16=41",
    "80-75=5",
    "60-33=27",
    "28+27=55",
    "65+34=99",
    "85+5=90",
    "26+59=85",
    "73-67=6",
    "33-23=10",
    "63+25=88",
    "57-22=35",
    "19+43=62",
    "66-31=35",
    "2+21=23",
    "57-38=19",
    "18+26=44",
    "58+22=80",
    "59+39=98",
    "89-50=39",
    "63-39=24",
    "88+2=90",
    "49-42=7",
    "50+24=74",
    "60-45=15",
    "59+27=86",
    "79-43=36",
    "98-69=29",
    "66-22=44",
    "81+15=96",
    "32+66=98",
    "33+17=50"
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Updated cells:" $idx
